$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host ($wb | Get-Member | Out-String)
